# Generate Report for Handoff
# Updates status of "7a03f314-a439-40fb-95f8-a2bd208697c9.md" to "Ready for handoff"
# on the Overview sheet and the per-locale sheets, and refreshes the
# "Latest Handoff Datetime" for every file that is now ready for / already
# handed off (i.e. everything except files still "In Translation").

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-03-08 07:12:06"
$wsZh.Range("D3").Value = "2016-03-08 07:12:06"
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "2016-03-08 07:12:06"
$wsZh.Range("D6").Value = "2016-03-08 07:12:06"
$wsZh.Range("D7").Value = "2016-03-08 07:12:06"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-03-08 07:12:17"
$wsDe.Range("D3").Value = "2016-03-08 07:12:17"
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "2016-03-08 07:12:17"
$wsDe.Range("D6").Value = "2016-03-08 07:12:17"
$wsDe.Range("D7").Value = "2016-03-08 07:12:17"
